# Auto-generated edit script applying the committed data refresh
# to the 北京-漫展信息.xlsx workbook (gh-pages data snapshot update).
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws0 = $wb.Worksheets.Item("展览")
$ws0.Range("F2").Value = 34
$ws0.Range("F4").Value = 26
$ws0.Range("F5").Value = 8820
$ws0.Range("F6").Value = 0
$ws0.Range("F7").Value = 234
$ws0.Range("F8").Value = 7092
$ws0.Range("F10").Value = 5369
$ws0.Range("F11").Value = 10
$ws0.Range("F12").Value = 6064
$ws0.Range("F13").Value = 1094
$ws0.Range("F14").Value = 389
$ws0.Range("F15").Value = 0
$ws0.Range("F17").Value = 0
$ws0.Range("F19").Value = 0
$ws0.Range("F20").Value = 0
$ws0.Range("F21").Value = 0
$ws0.Range("F22").Value = 0
$ws0.Range("F23").Value = 103
$ws0.Range("F25").Value = 1875
$ws0.Range("F26").Value = 1735
$ws0.Range("F29").Value = 2022
$ws0.Range("F30").Value = 78
$ws0.Range("F32").Value = 157
$ws0.Range("F33").Value = 1028
$ws0.Range("F34").Value = 2022
$ws0.Range("F35").Value = 303
$ws0.Range("F36").Value = 1350
$ws0.Range("F38").Value = 5143
$ws0.Range("D39").Value = "金蝉西路甲1号 酷车小镇"
$ws0.Range("E39").Value = "2024.08.24 10:00-08.25 17:00"
$ws0.Range("F39").Value = 0
$ws0.Range("I39").Value = "//i2.hdslb.com/bfs/openplatform/202407/3EF1Am6T1720430616435.jpeg"
$ws0.Range("F41").Value = 0
$ws0.Range("F42").Value = 0
$ws0.Range("F43").Value = 0
$ws0.Range("F45").Value = 1072
$ws0.Range("F46").Value = 971
$ws0.Range("F47").Value = 1348
$ws0.Range("F49").Value = 1091

# --- Sheet: 演出 ---
$ws1 = $wb.Worksheets.Item("演出")
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 5
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 31
$ws1.Range("F8").Value = 7
$ws1.Range("F9").Value = 40
$ws1.Range("F10").Value = 0
$ws1.Range("F12").Value = 9
$ws1.Range("F17").Value = 0

# --- Sheet: 全部类型 ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 34
$ws2.Range("F8").Value = 0
$ws2.Range("F12").Value = 3
$ws2.Range("F15").Value = 6064
$ws2.Range("F17").Value = 0
$ws2.Range("F18").Value = 396
$ws2.Range("F19").Value = 0
$ws2.Range("F20").Value = 552
$ws2.Range("F21").Value = 326
$ws2.Range("F22").Value = 269
$ws2.Range("F23").Value = 144
$ws2.Range("F24").Value = 202
$ws2.Range("F25").Value = 161
$ws2.Range("F26").Value = 0
$ws2.Range("F27").Value = 0
$ws2.Range("F28").Value = 9928
$ws2.Range("F29").Value = 1875
$ws2.Range("F30").Value = 0
$ws2.Range("F32").Value = 0
$ws2.Range("F33").Value = 78
$ws2.Range("F34").Value = 0
$ws2.Range("F35").Value = 0
$ws2.Range("F37").Value = 2022
$ws2.Range("F38").Value = 303
$ws2.Range("F39").Value = 1350
$ws2.Range("F40").Value = 0
$ws2.Range("F41").Value = 1193
$ws2.Range("F43").Value = 101
$ws2.Range("F44").Value = 167
$ws2.Range("F45").Value = 0
$ws2.Range("F46").Value = 0
$ws2.Range("F50").Value = 1091
